# Update the "Förändrad" (changed) date column (C) for all data rows
# from 45204 (2023-10-05) to 45205 (2023-10-06).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
